$d = $word.ActiveDocument

# 1) Footer paragraph: split into three runs
$d.Content.Find.Execute(
    "En el footer encontramos únicamente unos accesos directos a las redes sociales del club tragamillas. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "En el footer encontramos únicamente unos accesos directos a las redes sociales del club tragamillas junto a los enlaces disponibles en el menú. ",
    2
)

# 2) "y rojo son los colores..." -> "y de rojo son los colores..."
$d.Content.Find.Execute(
    "y rojo son los colores por defecto de esas tonalidades, no es ninguna variación.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "y de rojo son los colores por defecto de esas tonalidades, no es ninguna variación.",
    2
)

# 3) "la pagina del club tragamillas. " -> "la pagina del club tragamillas." (trailing space removed)
$d.Content.Find.Execute(
    "la pagina del club tragamillas. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "la pagina del club tragamillas.",
    2
)
